# Edit /tmp/work/before.xlsx to match the target revision.
#
# Semantics of the change (derived from the OOXML diff):
#  - Row 13 (an "orphan" row holding only the docente B/C text, with no
#    label in column A) is deleted entirely, shifting rows 14-22 up to
#    become rows 13-21 (heights/styles travel with the rows).
#  - After that shift, a handful of B/C (and mirrored C) cells get new
#    text so the sheet ends up with the content actually shown in the
#    target file (the diff shows several long paragraphs swapped for much
#    shorter ones, plus the "8452037 - Elisabeth ..." docente line moving
#    up into the "Objetivos:" row and again into the "Método:" row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the orphan row 13 - everything below shifts up by one row.
$ws.Rows.Item(13).Delete()

# 2) Apply the content updates on top of the shifted layout.

# Row 10 ("Objetivos:") - was the long PT objectives paragraph, now holds
# the docente identification line.
$ws.Range("B10:C10").Value = "8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara"

# Row 13 ("Programa resumido:") - was the PT short-syllabus bullet list,
# now just "Semestral".
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 ("Programa:") - was the long PT programa paragraph, now just the
# activation date string.
$ws.Range("B15:C15").Value = "01/01/2017"

# Row 18 ("Método:") - now also holds the docente identification line.
$ws.Range("B18:C18").Value = "8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara"

# Row 19 ("Critério:") - now holds the "A cada semestre ..." paragraph
# (previously under "Método:").
$ws.Range("B19:C19").Value = "A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas – a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo."

# Row 20 ("Norma de recuperação:") - now holds the "Sendo uma atividade
# prática ..." paragraph (previously under "Critério:").
$ws.Range("B20:C20").Value = "Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical."

# Row 21 ("Bibliografia:") - now holds "não tem" (previously under "Norma
# de recuperação:"); the old bibliography paragraph is dropped entirely.
$ws.Range("B21:C21").Value = "não tem"
